# Updates the cryptos list (Price / Volume(1h) columns) to new scraped values.
# Column D = Price, Column E = Volume(1h) (text values with leading/trailing
# double-space padding and a trailing "  " after the percent sign).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (column D, only when changed), new Volume (column E, only when changed)
$updates = @(
    @{ Row = 2;  D = "25.957.20";   E = "  -0.79%  " },
    @{ Row = 3;  D = "1.638.49";    E = "  -0.05%  " },
    @{ Row = 4;  E = "  +1.01%  " },
    @{ Row = 5;  D = "214.80";      E = "  -0.46%  " },
    @{ Row = 6;  E = "  +0.49%  " },
    @{ Row = 7;  E = "  +0.90%  " },
    @{ Row = 8;  E = "  -0.89%  " },
    @{ Row = 9;  D = "0.0637";      E = "  +0.37%  " },
    @{ Row = 10; D = "19.65";       E = "  -0.97%  " },
    @{ Row = 11; D = "0.0794";      E = "  +0.82%  " },
    @{ Row = 12; D = "1.864.55" },
    @{ Row = 13; D = "4.26";        E = "  -0.02%  " },
    @{ Row = 14; D = "1.625.32";    E = "  -0.62%  " },
    @{ Row = 15; E = "  -1.73%  " },
    @{ Row = 16; D = "0.0$([char]0x2083)0757" },
    @{ Row = 17; D = "62.55";       E = "  -1.25%  " },
    @{ Row = 18; D = "25.965.67";   E = "  -0.62%  " },
    @{ Row = 19; E = "  +0.95%  " },
    @{ Row = 20; D = "193.94";      E = "  +0.21%  " },
    @{ Row = 21; E = "  -1.98%  " },
    @{ Row = 22; E = "  -0.66%  " },
    @{ Row = 23; E = "  -1.60%  " },
    @{ Row = 24; D = "144.26";      E = "  +1.66%  " },
    @{ Row = 25; E = "  +0.13%  " },
    @{ Row = 26; E = "  +0.93%  " },
    @{ Row = 27; E = "  +2.98%  " },
    @{ Row = 28; E = "  -0.66%  " },
    @{ Row = 29; D = "15.50";       E = "  -0.74%  " },
    @{ Row = 30; E = "  -0.47%  " },
    @{ Row = 31; E = "  +0.71%  " },
    @{ Row = 32; E = "  -1.12%  " },
    @{ Row = 33; E = "  -0.33%  " },
    @{ Row = 34; D = "1.55";        E = "  -2.96%  " },
    @{ Row = 35; E = "  +1.86%  " },
    @{ Row = 36; D = "0.904";       E = "  -0.47%  " },
    @{ Row = 37; D = "1.139.17";    E = "  -0.46%  " },
    @{ Row = 38; E = "  -0.24%  " },
    @{ Row = 39; E = "  -1.56%  " },
    @{ Row = 40; E = "  +0.48%  " },
    @{ Row = 41; D = "99.35";       E = "  -0.86%  " },
    @{ Row = 42; D = "0.801";       E = "  +1.05%  " },
    @{ Row = 43; E = "  -2.78%  " },
    @{ Row = 44; D = "1.774.45";    E = "  -0.12%  " },
    @{ Row = 45; E = "  +7.25%  " },
    @{ Row = 46; E = "  +1.08%  " },
    @{ Row = 47; E = "  +2.79%  " },
    @{ Row = 48; E = "  +0.01%  " },
    @{ Row = 49; E = "  +0.04%  " },
    @{ Row = 50; D = "7.62";        E = "  +0.12%  " },
    @{ Row = 51; D = "0.0962";      E = "  -1.17%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Force text format so numeric-looking strings (e.g. "214.80",
        # "0.0637") keep their exact textual representation instead of
        # being coerced into a floating point number by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
